$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns B..G right into C..H (making room for the new
# "guest_gender" column at B), without touching the <cols> width definitions.
for ($col = 7; $col -ge 2; $col--) {
    $srcHeader = $ws.Cells.Item(1, $col)
    $srcData   = $ws.Cells.Item(2, $col)
    $dstHeader = $ws.Cells.Item(1, $col + 1)
    $dstData   = $ws.Cells.Item(2, $col + 1)

    $dstHeader.Value = $srcHeader.Value2

    $dstData.Style = "Normal"
    if ($srcData.NumberFormat -ne "General") {
        $dstData.NumberFormat = $srcData.NumberFormat
    }
    $dstData.Value = $srcData.Value2
}

# New "guest_gender" column in column B
$ws.Cells.Item(1, 2).Value = "guest_gender"
$ws.Cells.Item(2, 2).Value = "Male"

# Update selection to match the target workbook state
$ws.Range("G5").Select() | Out-Null
